# Update Input sheet values: Elevation Angle, Estimated GS Losses,
# Estimated S/C Losses, Atmospheric Losses - found the minimum angle
# needed to close the S-band link.
$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("Input")
$wsInput.Range("C6").Value = 20
$wsInput.Range("C16").Value = 2
$wsInput.Range("C20").Value = 2.2000000000000002
$wsInput.Range("C24").Value = 1.1000000000000001

# Update the selection left on the Input sheet.
$wsInput.Activate()
$wsInput.Range("C7").Select()

# Leave the selection on the UHF sheet as-is (still A10), it's just no
# longer the active tab once we move on to S-Band below.
$wsUHF = $wb.Worksheets.Item("UHF")
$wsUHF.Activate()
$wsUHF.Range("A10").Select()

# Make the S-Band sheet the active tab (this is now where the work is
# happening - finding the minimum elevation angle to close the link),
# scroll down a bit and reposition the selection.
$wsSBand = $wb.Worksheets.Item("S-Band")
$wsSBand.Activate()
$excel.ActiveWindow.ScrollRow = 13
$wsSBand.Range("A45").Select()
